# GOST_DbLib.xlsx edit script
# Commit: "Added: TO-220 RLB JRB"
#   - New "Analog ICs" sheet (NE555DR)
#   - New "RLB Inductors" sheet (RLB0914-102KL inductor)
#   - New "El Cap" sheet (ECAP (K50-35) electrolytic capacitor)
#   - New diode row (1N5404) on the "Diodes" sheet + header rename Value->Comment
#   - New MOSFET row (STP11NK50Z, TO-220) on the "Transistors" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Diodes sheet: rename column B header, add new diode part row
# ---------------------------------------------------------------------------
$diodes = $wb.Worksheets.Item("Diodes")
$diodes.Range("B1").Value = "Comment"

$diodes.Range("A4").Value = "1N5404"
$diodes.Range("B4").Value = "1N5404"
$diodes.Range("C4").Value = "Diode-CA"
$diodes.Range("D4").Value = "DO-201AD"
$diodes.Range("J4").Value = "Diodes.SchLib"
$diodes.Range("K4").Value = "DO Package.PcbLib"

# ---------------------------------------------------------------------------
# 2. Transistors sheet: add new TO-220 MOSFET row
# ---------------------------------------------------------------------------
$transistors = $wb.Worksheets.Item("Transistors")

$transistors.Range("A8").Value = "STP11NK50Z"
$transistors.Range("B8").Value = "MOSFET-N GDS"
$transistors.Range("C8").Value = "TO-220"
$transistors.Range("I8").Value = "Transistors.SchLib"
$transistors.Range("J8").Value = "TO Package.PcbLib"
$transistors.Range("K8").Value = "500В"
$transistors.Range("L8").Value = "10А"
$transistors.Range("M8").Value = "0,48Ом"
$transistors.Range("Q8").Value = "STP11NK50Z"

# ---------------------------------------------------------------------------
# 3. New sheet: Analog ICs (after OpAmps)
# ---------------------------------------------------------------------------
$opAmps = $wb.Worksheets.Item("OpAmps")
$analogIcs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $opAmps)
$analogIcs.Name = "Analog ICs"

$opAmps.Range("A1:F1").Copy()
$analogIcs.Range("A1:F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$analogIcs.Range("A1").Value = "Part Number"
$analogIcs.Range("B1").Value = "Library Ref"
$analogIcs.Range("C1").Value = "Footprint Ref"
$analogIcs.Range("D1").Value = "Library Path"
$analogIcs.Range("E1").Value = "Footprint Path"
$analogIcs.Range("F1").Value = "Comment"

$analogIcs.Range("A2").Value = "NE555DR"
$analogIcs.Range("B2").Value = "NE555DR"
$analogIcs.Range("C2").Value = "SO-8"
$analogIcs.Range("D2").Value = "Analog Ics.SchLib"
$analogIcs.Range("E2").Value = "SO Package.PcbLib"
$analogIcs.Range("F2").Value = "NE555DR"

# ---------------------------------------------------------------------------
# 4. New sheet: RLB Inductors (after Analog ICs)
# ---------------------------------------------------------------------------
$rlb = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $analogIcs)
$rlb.Name = "RLB Inductors"

$opAmps.Range("A1:F1").Copy()
$rlb.Range("A1:F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$rlb.Range("A1").Value = "Part Number"
$rlb.Range("B1").Value = "Library Ref"
$rlb.Range("C1").Value = "Footprint Ref"
$rlb.Range("D1").Value = "Library Path"
$rlb.Range("E1").Value = "Footprint Path"
$rlb.Range("F1").Value = "Comment"
$rlb.Range("G1").Value = "Inductance"

$rlb.Range("A2").Value = "RLB0914-102KL"
$rlb.Range("B2").Value = "Core Inductor"
$rlb.Range("C2").Value = "RLB0914"
$rlb.Range("D2").Value = "Inductors.SchLib"
$rlb.Range("E2").Value = "RLB.PcbLib"
$rlb.Range("F2").Value = "RLB0914-102KL"
$rlb.Range("G2").Value = "1000мк"

$rlb.PageSetup.PaperSize = 9
$rlb.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. New sheet: El Cap (after RLB Inductors)
# ---------------------------------------------------------------------------
$elcap = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $rlb)
$elcap.Name = "El Cap"

$opAmps.Range("A1:F1").Copy()
$elcap.Range("A1:F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$elcap.Range("A1").Value = "Part Number"
$elcap.Range("B1").Value = "Library Ref"
$elcap.Range("C1").Value = "Footprint Ref"
$elcap.Range("D1").Value = "Library Path"
$elcap.Range("E1").Value = "Footprint Path"
$elcap.Range("F1").Value = "Comment"
$elcap.Range("G1").Value = "Value"
$elcap.Range("H1").Value = "Voltage"

$elcap.Range("A2").Value = "ECAP (К50-35)"
$elcap.Range("B2").Value = "Polarized Capacitor"
$elcap.Range("C2").Value = "JRB12.5/13"
$elcap.Range("D2").Value = "Capacitors.SchLib"
$elcap.Range("E2").Value = "JRB.PcbLib"
$elcap.Range("F2").Value = "ECAP (К50-35)"
$elcap.Range("G2").Value = "22мк"
$elcap.Range("H2").Value = "350В"

# ---------------------------------------------------------------------------
# 6. Sheet view / selection bookkeeping
#    (the last-activated sheet/range ends up as the active tab & selection)
# ---------------------------------------------------------------------------
$diodes.Range("C10").Select()

$opAmps.Rows.Item(1).Select()

$transistors.Range("L8").Select()
